# Generate Report for Handoff
#
# Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
# for the d57fe60c-19d4-43bc-856f-07ba113e1615.md row across the
# Overview / zh-cn / de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-09-05 18:47:14"

# --- zh-cn sheet: "Latest Handoff Datetime" column (H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-09-05 18:47:08"

# --- de-de sheet: "Latest Handoff Datetime" column (H) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-09-05 18:47:14"
